# Updates the "Updated symbol list" coin-ranking snapshot sheet.
# Values are numeric-looking text (prices, percentages) stored as TEXT in the
# original workbook, so each write forces a text NumberFormat first (otherwise
# COM auto-coerces strings like "302.07" into real numbers), then restores the
# "Normal" cell style so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$cellRef, [string]$value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '302.07'
Set-TextCell 'E2' '-4.17%'
# Row 3
Set-TextCell 'E3' '-1.91%'
# Row 4
Set-TextCell 'D4' '5.044'
Set-TextCell 'E4' '-1.77%'
# Row 5
Set-TextCell 'D5' '0.07991'
Set-TextCell 'E5' '-1.39%'
# Row 6
Set-TextCell 'D6' '1.935'
Set-TextCell 'E6' '-8.85%'
# Row 7
Set-TextCell 'B7' 'GateToken'
Set-TextCell 'C7' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell 'D7' '4.044'
Set-TextCell 'E7' '-2.49%'
# Row 8
Set-TextCell 'B8' 'KuCoinToken'
Set-TextCell 'C8' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextCell 'D8' '7.741'
Set-TextCell 'E8' '-3.36%'
# Row 9
Set-TextCell 'B9' 'BTSEToken'
Set-TextCell 'C9' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell 'D9' '2.949'
Set-TextCell 'E9' '9.06%'
# Row 10
Set-TextCell 'B10' 'MXToken'
Set-TextCell 'C10' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D10' '0.9218'
Set-TextCell 'E10' '-0.98%'
# Row 11
Set-TextCell 'B11' 'LiechtensteinCryptoassetsExchange'
Set-TextCell 'C11' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 'D11' '0.1296'
Set-TextCell 'E11' '28.39%'
# Row 12
Set-TextCell 'B12' 'WazirX'
Set-TextCell 'C12' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 'D12' '0.1849'
Set-TextCell 'E12' '-1.19%'
# Row 13
Set-TextCell 'B13' 'MandalaExchangeToken'
Set-TextCell 'C13' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 'D13' '0.09652'
Set-TextCell 'E13' '5.24%'
# Row 14
Set-TextCell 'B14' 'BitrueCoin'
Set-TextCell 'C14' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 'D14' '0.03625'
Set-TextCell 'E14' '0.27%'
# Row 15
Set-TextCell 'B15' 'BitMartToken'
Set-TextCell 'C15' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 'D15' '0.09854'
Set-TextCell 'E15' '-0.58%'
# Row 16
Set-TextCell 'B16' 'BitForexToken'
Set-TextCell 'C16' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 'D16' '0.001393'
Set-TextCell 'E16' '-4.16%'
# Row 17
Set-TextCell 'B17' 'TigerCash'
Set-TextCell 'C17' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell 'D17' '0.005810'
Set-TextCell 'E17' '1.64%'
# Row 18
Set-TextCell 'B18' 'LEO'
Set-TextCell 'C18' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 'D18' '3.505'
Set-TextCell 'E18' '1.02%'
# Row 19
Set-TextCell 'D19' '0.3431'
Set-TextCell 'E19' '1.82%'
# Row 20
Set-TextCell 'D20' '0.1310'
Set-TextCell 'E20' '-1.77%'
# Row 21
Set-TextCell 'D21' '5.049'
Set-TextCell 'E21' '-1.91%'
# Row 22
Set-TextCell 'E22' '9.19%'
# Row 23
Set-TextCell 'D23' '0.04524'
Set-TextCell 'E23' '-1.26%'
# Row 24
Set-TextCell 'D24' '0.001215'
Set-TextCell 'E24' '-2.86%'
# Row 25
Set-TextCell 'D25' '0.004819'
Set-TextCell 'E25' '2.46%'
# Row 26
Set-TextCell 'E26' '-0.34%'
# Row 27
Set-TextCell 'D27' '0.0003005'
Set-TextCell 'E27' '-33.50%'
# Row 39
Set-TextCell 'D39' '0.01905'
Set-TextCell 'E39' '-3.45%'
# Row 40
Set-TextCell 'D40' '0.04704'
Set-TextCell 'E40' '-3.98%'
# Row 41
Set-TextCell 'D41' '0.007518'
Set-TextCell 'E41' '-4.26%'
# Row 42
Set-TextCell 'D42' '0.009650'
Set-TextCell 'E42' '22.85%'
# Row 43
Set-TextCell 'D43' '0.1328'
Set-TextCell 'E43' '-4.55%'
# Row 44
Set-TextCell 'D44' '0.002111'
Set-TextCell 'E44' '-0.02%'
# Row 45
Set-TextCell 'D45' '0.01084'
Set-TextCell 'E45' '-6.64%'
# Row 46
Set-TextCell 'E46' '-4.51%'
# Row 47
Set-TextCell 'D47' '0.00000000750'
Set-TextCell 'E47' '-0.36%'
# Row 48
Set-TextCell 'E48' '82.08%'
# Row 49
Set-TextCell 'E49' '-21.93%'
# Row 50
Set-TextCell 'D50' '0.00002101'
Set-TextCell 'E50' '-0.36%'
# Row 51
Set-TextCell 'D51' '0.0002001'
Set-TextCell 'E51' '-0.36%'

Write-Host "Applied 96 cell updates"
